$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.02732947460379011
$ws.Range("B2").Value = 0.008099341481096272
$ws.Range("C2").Value = 0.04823371022939682
$ws.Range("D2").Value = 0.08814053599658071
$ws.Range("E2").Value = 0.1009959044756205
$ws.Range("F2").Value = 0.09413133853757837
$ws.Range("G2").Value = 0.02603567827558224

$ws.Range("A3").Value = 0.02100294696079608
$ws.Range("B3").Value = 0.006438110103244445
$ws.Range("C3").Value = -0.01669603289021728
$ws.Range("D3").Value = 0.05868442345439766
$ws.Range("E3").Value = 0.09840140495944422
$ws.Range("F3").Value = 0.07352196917387679
$ws.Range("G3").Value = 0.01758599876289391

$ws.Range("A4").Value = 0.05426793787707845
$ws.Range("B4").Value = 0.02507370232752118
$ws.Range("C4").Value = 0.06021424010396004
$ws.Range("D4").Value = 0.1191507621515193
$ws.Range("E4").Value = 0.1266817804555617
$ws.Range("F4").Value = 0.1228009158748118
$ws.Range("G4").Value = 0.02242512085877799

$ws.Range("A5").Value = 0.01324707873219017
$ws.Range("B5").Value = 0.004479855794709306
$ws.Range("C5").Value = -0.0004699623746217008
$ws.Range("D5").Value = 0.07868029157330532
$ws.Range("E5").Value = 0.09626063126655131
$ws.Range("F5").Value = 0.0865871108044384
$ws.Range("G5").Value = 0.01819391779573026
